$p = $ppt.ActivePresentation

# Slide 1: Title "Header" + " " + "with" + " " + [Courier "inline code"]
#          -> "Header with " (single run) + [Courier "inline code"] (unchanged)
$s1 = $p.Slides.Item(1)
$title1 = $s1.Shapes.Item(1).TextFrame.TextRange
$mid1 = $title1.Characters(7, 6)          # " with " between "Header" and "inline code"
[void]$mid1.Delete()
$head1 = $title1.Characters(1, 6)         # "Header"
[void]$head1.InsertAfter(" with ")

# Slide 2: Title "Syntax" + " " + "highlighting" -> "Syntax highlighting" (single run)
$s2 = $p.Slides.Item(2)
$title2 = $s2.Shapes.Item(1).TextFrame.TextRange
$title2.Text = "X"
$title2.Text = "Syntax highlighting"

# Slide 3: Title "Two" + " " + "column" + " " + "slide" -> "Two column slide" (single run)
$s3 = $p.Slides.Item(3)
$title3 = $s3.Shapes.Item(1).TextFrame.TextRange
$title3.Text = "X"
$title3.Text = "Two column slide"
